$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 713; this pushes the existing rows
# 713:743 down to 715:745 (matches the diff's row-shift + dimension growth
# to A1:R745) while leaving rows 1:712 untouched.
$ws.Rows.Item(713).Insert()
$ws.Rows.Item(713).Insert()

# Populate the two freshly inserted rows with the new weekly price record
# (week of 2023-08-09, serial 45147) for Ciboulette, "Primera" and
# "Segunda" quality grades.
$ws.Range("A713").Value = 6
$ws.Range("B713").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C713").Value = "Metropolitana"
$ws.Range("D713").Value = 45147
$ws.Range("E713").Value = 13
$ws.Range("F713").Value = 100112039
$ws.Range("G713").Value = "Ciboulette"
$ws.Range("H713").Value = "Sin especificar"
$ws.Range("I713").Value = "Primera"
$ws.Range("J713").Value = 680
$ws.Range("K713").Value = 1000
$ws.Range("L713").Value = 1200
$ws.Range("M713").Value = 1112
$ws.Range("N713").Value = "`$/docena de atados"
$ws.Range("O713").Value = "Región Metropolitana"
$ws.Range("P713").Value = 371
$ws.Range("Q713").Value = 3
$ws.Range("R713").Value = "Hortaliza"

$ws.Range("A714").Value = 6
$ws.Range("B714").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C714").Value = "Metropolitana"
$ws.Range("D714").Value = 45147
$ws.Range("E714").Value = 13
$ws.Range("F714").Value = 100112039
$ws.Range("G714").Value = "Ciboulette"
$ws.Range("H714").Value = "Sin especificar"
$ws.Range("I714").Value = "Segunda"
$ws.Range("J714").Value = 200
$ws.Range("K714").Value = 800
$ws.Range("L714").Value = 800
$ws.Range("M714").Value = 800
$ws.Range("N714").Value = "`$/docena de atados"
$ws.Range("O714").Value = "Región Metropolitana"
$ws.Range("P714").Value = 267
$ws.Range("Q714").Value = 3
$ws.Range("R714").Value = "Hortaliza"

# Make sure the date cells keep the date/time display format used by the
# rest of column D (style index 2 / numFmtId 165 in the original file).
$ws.Range("D713:D714").NumberFormat = "YYYY-MM-DD HH:MM:SS"
